# Update "想去人数" (interested-people count, column F) figures that were
# refreshed when the site's data was regenerated.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 8133
$ws1.Range("F5").Value  = 6505
$ws1.Range("F7").Value  = 2059
$ws1.Range("F8").Value  = 566
$ws1.Range("F15").Value = 8491
$ws1.Range("F20").Value = 1804
$ws1.Range("F25").Value = 29
$ws1.Range("F28").Value = 8
$ws1.Range("F30").Value = 2067
$ws1.Range("F31").Value = 844
$ws1.Range("F32").Value = 467
$ws1.Range("F35").Value = 174
$ws1.Range("F37").Value = 2
$ws1.Range("F38").Value = 22

# Sheet "本地生活" (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 305

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 8133
$ws4.Range("F7").Value  = 305
$ws4.Range("F9").Value  = 6505
$ws4.Range("F11").Value = 2059
$ws4.Range("F14").Value = 566
$ws4.Range("F23").Value = 8491
$ws4.Range("F27").Value = 1804
$ws4.Range("F31").Value = 8
$ws4.Range("F32").Value = 2067
$ws4.Range("F33").Value = 844
$ws4.Range("F35").Value = 467
